# Weekly data update: insert a new price record as row 146 on the Berenjena
# (Hortaliza, Vega Central Mapocho de Santiago) sheet. Inserting the row
# shifts the existing rows 146-224 down to 147-225, matching the target
# dimension A1:R225.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 146, pushing everything below it down
# by one row.
$ws.Rows(146).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(146, 1).Value  = 9
$ws.Cells.Item(146, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(146, 3).Value  = "Metropolitana"
$ws.Cells.Item(146, 4).Value  = 44572
$ws.Cells.Item(146, 5).Value  = 13
$ws.Cells.Item(146, 6).Value  = 100112001
$ws.Cells.Item(146, 7).Value  = "Berenjena"
$ws.Cells.Item(146, 8).Value  = "Sin especificar"
$ws.Cells.Item(146, 9).Value  = "Primera"
$ws.Cells.Item(146, 10).Value = 52
$ws.Cells.Item(146, 11).Value = 10000
$ws.Cells.Item(146, 12).Value = 12000
$ws.Cells.Item(146, 13).Value = 11000
$ws.Cells.Item(146, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(146, 15).Value = "Región Metropolitana"
$ws.Cells.Item(146, 16).Value = 183
$ws.Cells.Item(146, 17).Value = 60
$ws.Cells.Item(146, 18).Value = "Hortaliza"
